$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new task-log rows (Investigación / Angular) below the existing data.
$ws.Range("B17").Copy()
$ws.Range("B18").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B19").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A18").Value = "Federico Speroni"
$ws.Range("B18").Value = (Get-Date -Year 2017 -Month 4 -Day 19 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = "Investigación"
$ws.Range("E18").Value = "Angular"

$ws.Range("A19").Value = "Federico Speroni"
$ws.Range("B19").Value = (Get-Date -Year 2017 -Month 4 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = "Investigación"
$ws.Range("E19").Value = "Angular"

$ws.Range("E18").Select()
